$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row 80 first so the new shared strings are allocated in the
# same order as in the target file (name, then description), before the
# edited B31 description reuses a later slot.
$ws.Range("A80").Value = "listMetRxnsWithFluxes.m"
$ws.Range("B80").Value = "This code was copied from another project - TME modeling, and works. No specific testing was deemed needed, the output looks very reasonable."

# Update existing cell B31: test description text changed
$ws.Range("B31").Value = "This code is fairly trivial, just calls other functions and merges data, as well as a simple task analysis. No further testing was deemed needed"

# Update selection to match target (engine does not persist scroll
# position / topLeftCell, only the active cell/selection).
$ws.Activate()
$ws.Range("B31").Select()
